# Auto-generated COM-interop script to update cryptos worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) keeps its text formatting so values like "234.34"
# are not auto-converted to numbers by Excel.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.494.37'
$ws.Range("E2").Value = '  -0.09%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.080.66'
$ws.Range("E3").Value = '  +0.06%  '

# Row 4
$ws.Range("E4").Value = '  +0.09%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.34'
$ws.Range("E5").Value = '  -0.46%  '

# Row 6
$ws.Range("E6").Value = '  +1.23%  '

# Row 7
$ws.Range("E7").Value = '  +0.02%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '57.68'
$ws.Range("E8").Value = '  -1.37%  '

# Row 9
$ws.Range("E9").Value = '  +0.85%  '

# Row 10
$ws.Range("E10").Value = '  +2.06%  '

# Row 11
$ws.Range("E11").Value = '  +1.15%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.389.04'
$ws.Range("E12").Value = '  +0.50%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.42'
$ws.Range("E13").Value = '  -1.29%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.84'
$ws.Range("E14").Value = '  -1.66%  '

# Row 15
$ws.Range("E15").Value = '  +0.06%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.21'
$ws.Range("E16").Value = '  +0.00%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.094.83'
$ws.Range("E17").Value = '  +1.21%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.454.89'
$ws.Range("E18").Value = '  -0.65%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.19'
$ws.Range("E19").Value = '  -1.52%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.67'
$ws.Range("E20").Value = '  -0.90%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0820'
$ws.Range("E21").Value = '  +0.29%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '226.48'
$ws.Range("E22").Value = '  -0.36%  '

# Row 23
$ws.Range("E23").Value = '  -0.01%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.45'
$ws.Range("E24").Value = '  +2.45%  '

# Row 25
$ws.Range("E25").Value = '  -2.10%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '168.64'
$ws.Range("E26").Value = '  +1.27%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.90'
$ws.Range("E27").Value = '  -0.57%  '

# Row 28
$ws.Range("B28").Value = 'ImmutableX'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.44'
$ws.Range("E28").Value = '  -4.65%  '

# Row 29
$ws.Range("B29").Value = 'Kaspa'
$ws.Range("C29").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.131'
$ws.Range("E29").Value = '  +2.58%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.18'
$ws.Range("E30").Value = '  -0.94%  '

# Row 31
$ws.Range("E31").Value = '  -0.38%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.64'
$ws.Range("E32").Value = '  +1.76%  '

# Row 33
$ws.Range("E33").Value = '  -1.19%  '

# Row 34
$ws.Range("E34").Value = '  +0.26%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.55'
$ws.Range("E35").Value = '  -1.38%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.42'
$ws.Range("E36").Value = '  +1.65%  '

# Row 37
$ws.Range("E37").Value = '  +0.52%  '

# Row 38
$ws.Range("E38").Value = '  +0.19%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.59'
$ws.Range("E39").Value = '  -5.07%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.95'
$ws.Range("E40").Value = '  -0.37%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.490.08'
$ws.Range("E41").Value = '  +2.09%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0956'
$ws.Range("E42").Value = '  +0.22%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '97.15'
$ws.Range("E43").Value = '  +1.02%  '

# Row 44
$ws.Range("E44").Value = '  +0.11%  '

# Row 45
$ws.Range("E45").Value = '  -1.92%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.19'
$ws.Range("E46").Value = '  -9.24%  '

# Row 47
$ws.Range("E47").Value = '  +0.21%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '15.50'
$ws.Range("E48").Value = '  -2.50%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.27'
$ws.Range("E49").Value = '  +0.05%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.98'
$ws.Range("E50").Value = '  +1.33%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.275.80'
$ws.Range("E51").Value = '  +0.52%  '
